# Updated cryptos list (Price + Volume(1h)) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.853.68"
$ws.Range("E2").Value = "  -1.02%  "
# Row 3
$ws.Range("D3").Value = "1.856.76"
$ws.Range("E3").Value = "  -0.66%  "
# Row 4
$ws.Range("E4").Value = "  -0.17%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.81"
# Row 6
$ws.Range("E6").Value = "  -0.13%  "
# Row 7
$ws.Range("E7").Value = "  -1.54%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3657"
$ws.Range("E8").Value = "  -2.73%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07173"
$ws.Range("E9").Value = "  +0.11%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8895"
$ws.Range("E10").Value = "  +0.51%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.68"
$ws.Range("E11").Value = "  -0.58%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07534"
$ws.Range("E12").Value = "  -0.41%  "
# Row 13
$ws.Range("D13").Value = "1.846.37"
$ws.Range("E13").Value = "  -1.26%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "91.67"
$ws.Range("E14").Value = "  +2.67%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.229"
$ws.Range("E15").Value = "  -1.89%  "
# Row 16
$ws.Range("E16").Value = "  -0.17%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008530"
$ws.Range("E17").Value = "  -0.26%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.05"
$ws.Range("E18").Value = "  -0.94%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("E19").Value = "  -0.13%  "
# Row 20
$ws.Range("D20").Value = "26.894.77"
$ws.Range("E20").Value = "  -1.08%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.020"
$ws.Range("E21").Value = "  -0.26%  "
# Row 22
$ws.Range("D22").Value = "2.085.45"
$ws.Range("E22").Value = "  -1.97%  "
# Row 23
$ws.Range("E23").Value = "  -2.97%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.445"
$ws.Range("E24").Value = "  -0.54%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.09"
$ws.Range("E25").Value = "  -3.58%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.805"
$ws.Range("E26").Value = "  -2.56%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.83"
$ws.Range("E27").Value = "  -1.18%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.054"
$ws.Range("E28").Value = "  -5.24%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.98"
$ws.Range("E29").Value = "  +0.02%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.638"
$ws.Range("E30").Value = "  -2.22%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.665"
$ws.Range("E31").Value = "  -0.63%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09219"
$ws.Range("E32").Value = "  +2.15%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05095"
$ws.Range("E33").Value = "  -1.30%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.020"
$ws.Range("E34").Value = "  -2.41%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7341"
$ws.Range("E35").Value = "  -2.97%  "
# Row 36
$ws.Range("E36").Value = "  -2.25%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.220"
$ws.Range("E37").Value = "  +6.07%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02006"
$ws.Range("E38").Value = "  -1.69%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.474"
$ws.Range("E39").Value = "  -1.53%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.073"
$ws.Range("E40").Value = "  -0.91%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5319"
$ws.Range("E41").Value = "  -1.66%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "118.45"
$ws.Range("E42").Value = "  +3.15%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.497"
$ws.Range("E43").Value = "  -2.59%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.361"
$ws.Range("E44").Value = "  -2.04%  "
# Row 45
$ws.Range("E45").Value = "  -1.22%  "
# Row 46
$ws.Range("E46").Value = "  -1.11%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9995"
$ws.Range("E47").Value = "  -0.16%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.961"
$ws.Range("E48").Value = "  -1.96%  "
# Row 49
$ws.Range("E49").Value = "  -1.19%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.96"
$ws.Range("E50").Value = "  +1.22%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.84"
$ws.Range("E51").Value = "  -3.34%  "
